$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D:E").Insert(-4161, 0)

# Copy cell formatting (number formats/styles) from column F (the old column D, now shifted)
# into the new D:E columns so the new quarter columns match formatting of existing data.
# Only the row blocks that actually contain quarterly data are touched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Populate the two new quarter columns (D = Dec-2018, E = Sep-2018) with their reported values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 611900
$ws.Range("E8").Value = 634000
$ws.Range("D9").Value = 446000
$ws.Range("E9").Value = 468900
$ws.Range("D10").Value = 165900
$ws.Range("E10").Value = 165100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 53400
$ws.Range("E15").Value = 53600
$ws.Range("D17").Value = 562000
$ws.Range("E17").Value = 594300
$ws.Range("D18").Value = 49900
$ws.Range("E18").Value = 39700
$ws.Range("D20").Value = -4400
$ws.Range("E20").Value = 1400
$ws.Range("D21").Value = 98900
$ws.Range("E21").Value = 94700
$ws.Range("D22").Value = 36000
$ws.Range("E22").Value = 35300
$ws.Range("D23").Value = 9500
$ws.Range("E23").Value = 5800
$ws.Range("D24").Value = -4400
$ws.Range("E24").Value = -2700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 13900
$ws.Range("E26").Value = 8500
$ws.Range("D27").Value = 5000
$ws.Range("E27").Value = 2100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -10900
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4400
$ws.Range("E32").Value = -1400
$ws.Range("D33").Value = -5900
$ws.Range("E33").Value = 2100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -5900
$ws.Range("E35").Value = 2100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 86700
$ws.Range("E41").Value = 70000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 562500
$ws.Range("E43").Value = 662600
$ws.Range("D44").Value = 49400
$ws.Range("E44").Value = 49200
$ws.Range("D45").Value = 17300
$ws.Range("E45").Value = 19300
$ws.Range("D46").Value = 715800
$ws.Range("E46").Value = 801100
$ws.Range("D47").Value = 274000
$ws.Range("E47").Value = 277000
$ws.Range("D48").Value = 3457300
$ws.Range("E48").Value = 3450800
$ws.Range("D49").Value = 622300
$ws.Range("E49").Value = 630700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 140800
$ws.Range("E52").Value = 138200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5210300
$ws.Range("E54").Value = 5297800
$ws.Range("D57").Value = 498500
$ws.Range("E57").Value = 546300
$ws.Range("D58").Value = 6000
$ws.Range("E58").Value = 6000
$ws.Range("D59").Value = 132700
$ws.Range("E59").Value = 129500
$ws.Range("D60").Value = 637200
$ws.Range("E60").Value = 681900
$ws.Range("D61").Value = 2278800
$ws.Range("E61").Value = 2619500
$ws.Range("D62").Value = 94300
$ws.Range("E62").Value = 87100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3359800
$ws.Range("E66").Value = 3388500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 359700
$ws.Range("E70").Value = 353300
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -74000
$ws.Range("E72").Value = -74500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1490800
$ws.Range("E76").Value = 1556000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -5900
$ws.Range("E81").Value = 2100
$ws.Range("D83").Value = 53400
$ws.Range("E83").Value = 53600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 133300
$ws.Range("E89").Value = 40000
$ws.Range("D91").Value = -87300
$ws.Range("E91").Value = -69100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -81500
$ws.Range("E94").Value = -69900
$ws.Range("D96").Value = -37000
$ws.Range("E96").Value = -37000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -33500
$ws.Range("E100").Value = 43900
$ws.Range("D101").Value = -1600
$ws.Range("E101").Value = 800
$ws.Range("D102").Value = 16700
$ws.Range("E102").Value = 14700
